# Adds the "square 2" (column H) data, some custom parameter values in
# square 4 / square 5, and renames the terra_*.png border images to
# terra_*.jpg across all data rows ("Add middle squares. Add custom
# parameters").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 2 (card 001)
$ws.Range("H2").Value = "scissors.png"
$ws.Range("L2").Value = "terra_orange.jpg"

# Row 3 (card 002)
$ws.Range("H3").Value = "paper.png"
$ws.Range("J3").Value = 36
$ws.Range("L3").Value = "terra_green.jpg"

# Row 4 (card 003) -- square 1 (G) moves out, square 5 (K) gets it instead
$ws.Range("G4").Clear()
$ws.Range("K4").Value = "paper.png"
$ws.Range("L4").Value = "terra_blue.jpg"

# Row 5 (card 004) -- square 1 (G) moves out, square 5 (K) gets it instead
$ws.Range("G5").Clear()
$ws.Range("K5").Value = "rock.png"
$ws.Range("L5").Value = "terra_blue.jpg"

# Row 6 (card 005)
$ws.Range("L6").Value = "terra_purple.jpg"

# Row 7 (card 006)
$ws.Range("K7").Value = "scissors.png"
$ws.Range("L7").Value = "terra_orange.jpg"

# Row 8 (card 007)
$ws.Range("L8").Value = "terra_orange.jpg"

# Row 9 (card 008)
$ws.Range("L9").Value = "terra_purple.jpg"

# Row 10 (card 009) -- old numeric square 5 (K=2) is dropped, square 2 (H) added
$ws.Range("H10").Value = "rock.png"
$ws.Range("K10").Clear()
$ws.Range("L10").Value = "terra_orange.jpg"

# Row 11 (card 010)
$ws.Range("H11").Value = "paper.png"
$ws.Range("J11").Value = 20
$ws.Range("L11").Value = "terra_blue.jpg"

# Row 12 (card 011)
$ws.Range("L12").Value = "terra_orange.jpg"

# Row 13 (card 012)
$ws.Range("H13").Value = "ABC"
$ws.Range("J13").Value = 40
$ws.Range("K13").Value = "3Y"
$ws.Range("L13").Value = "terra_green.jpg"

# Row 14 (card 013)
$ws.Range("H14").Value = "XYZ"
$ws.Range("K14").Value = "4Z"
$ws.Range("L14").Value = "terra_green.jpg"

# ---------------------------------------------------------------------
# Column widths: column H now matches G's width (square 1/2 are grouped),
# and K/L pick up new best-fit widths to match the new content.
$ws.Columns.Item(8).ColumnWidth = 10.665
$ws.Columns.Item(11).ColumnWidth = 10.665
$ws.Columns.Item(12).ColumnWidth = 15

Write-Host "done"
